$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows (270-301), extending the series through 2021-06-28 ("aggiornamento fino a 28/06 incluso").
# Columns: A=date serial, B=nuovi pos., C=somma mobile 7gg., D=somma mobile 7gg. per 100mila abitanti
$data = @(
    @(44344, 0, 0, 0),
    @(44345, 0, 0, 0),
    @(44346, 0, 0, 0),
    @(44347, 0, 0, 0),
    @(44348, 0, 0, 0),
    @(44349, 0, 0, 0),
    @(44350, 0, 0, 0),
    @(44351, 0, 0, 0),
    @(44352, 0, 0, 0),
    @(44353, 1, 1, 18.93939393939394),
    @(44354, 0, 1, 18.93939393939394),
    @(44355, 0, 1, 18.93939393939394),
    @(44356, 0, 1, 18.93939393939394),
    @(44357, 0, 1, 18.93939393939394),
    @(44358, 0, 1, 18.93939393939394),
    @(44359, 0, 1, 18.93939393939394),
    @(44360, 0, 0, 0),
    @(44361, 0, 0, 0),
    @(44362, 0, 0, 0),
    @(44363, 0, 0, 0),
    @(44364, 1, 1, 18.93939393939394),
    @(44365, 0, 1, 18.93939393939394),
    @(44366, 0, 1, 18.93939393939394),
    @(44367, 1, 2, 37.87878787878788),
    @(44368, 0, 2, 37.87878787878788),
    @(44369, 0, 2, 37.87878787878788),
    @(44370, 1, 3, 56.81818181818181),
    @(44371, 1, 3, 56.81818181818181),
    @(44372, 0, 3, 56.81818181818181),
    @(44373, 1, 4, 75.75757575757575),
    @(44374, 1, 4, 75.75757575757575),
    @(44375, 1, 5, 94.6969696969697)
)

$startRow = 270
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $srcRange = $ws.Range("A269:D269")
    $dstRange = $ws.Range("A" + $row + ":D" + $row)
    # Copy formatting from the last existing row so the new row keeps the same date style (col A) and plain number style (cols B-D).
    $srcRange.Copy($dstRange)
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

$wb.Save()
